$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "34.089.71"
Set-TextValue "E2" "  +0.10%  "
Set-TextValue "D3" "1.780.58"
Set-TextValue "E3" "  -0.46%  "
Set-TextValue "E4" "  +0.29%  "
Set-TextValue "E5" "  -0.60%  "
Set-TextValue "E6" "  -0.04%  "
Set-TextValue "E7" "  +0.28%  "
Set-TextValue "D8" "31.82"
Set-TextValue "E8" "  -1.48%  "
Set-TextValue "D10" "0.0686"
Set-TextValue "E10" "  +0.33%  "
Set-TextValue "E11" "  +0.74%  "
Set-TextValue "D12" "2.036.76"
Set-TextValue "E12" "  -0.48%  "
Set-TextValue "B13" "Chainlink"
Set-TextValue "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D13" "10.92"
Set-TextValue "E13" "  -4.65%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.773.58"
Set-TextValue "E14" "  -0.81%  "
Set-TextValue "D15" "34.083.23"
Set-TextValue "E15" "  +0.17%  "
Set-TextValue "D16" "0.621"
Set-TextValue "E16" "  -0.13%  "
Set-TextValue "E17" "  -0.36%  "
Set-TextValue "D18" "67.61"
Set-TextValue "E18" "  -0.54%  "
Set-TextValue "D19" "244.75"
Set-TextValue "E20" "  +1.85%  "
Set-TextValue "E22" "  +1.50%  "
Set-TextValue "E23" "  -0.01%  "
Set-TextValue "E24" "  -1.04%  "
Set-TextValue "D25" "161.67"
Set-TextValue "E25" "  -0.29%  "
Set-TextValue "D26" "7.09"
Set-TextValue "E26" "  -0.69%  "
Set-TextValue "D27" "16.22"
Set-TextValue "E27" "  +0.01%  "
Set-TextValue "E28" "  +0.69%  "
Set-TextValue "E29" "  +0.34%  "
Set-TextValue "E31" "  +0.00%  "
Set-TextValue "D32" "3.69"
Set-TextValue "E32" "  +1.34%  "
Set-TextValue "E33" "  +2.18%  "
Set-TextValue "E34" "  -2.49%  "
Set-TextValue "D35" "1.445.71"
Set-TextValue "E35" "  +3.31%  "
Set-TextValue "D36" "2.46"
Set-TextValue "E36" "  +4.35%  "
Set-TextValue "E37" "  -0.60%  "
Set-TextValue "E38" "  +0.95%  "
Set-TextValue "E39" "  -0.63%  "
Set-TextValue "E40" "  +1.81%  "
Set-TextValue "D41" "80.28"
Set-TextValue "E41" "  +0.12%  "
Set-TextValue "E42" "  +1.24%  "
Set-TextValue "E43" "  -0.59%  "
Set-TextValue "D44" "13.66"
Set-TextValue "E44" "  -0.52%  "
Set-TextValue "E45" "  +2.45%  "
Set-TextValue "E46" "  +0.20%  "
Set-TextValue "E47" "  -0.97%  "
Set-TextValue "D48" "1.937.89"
Set-TextValue "E48" "  -0.48%  "
Set-TextValue "E49" "  +0.30%  "
Set-TextValue "D50" "104.16"
Set-TextValue "E50" "  -3.27%  "
Set-TextValue "E51" "  -6.77%  "
